$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "关键词/Provide a list of 20 most related keywords, in the following format:     - Keyword 1     - Keyword 2     - Keyword 3 ......"
$ws.Range("F3").Value = "大纲/提供20个最重要最相关的知识点大纲来全面深入学习主题相关内容, in the following format:     - sentence 1     - sentence 2     - sentence 3   ......"
$ws.Range("F4").Value = "专家发言/提供10个该领域最权威的专家的专业解答.`nFinal output are in the following format:     - 段落 1     - 段落 2     - 段落 3  ......`n"
$ws.Range("F5").Value = "Q&A/Provide a list of 20 most related best questions with answers, in this format:   {question}/{answer}.  Final output are in the following format:     - item 1     - item 2     - item 3 ......"
$ws.Range("F6").Value = "推荐书/Provide a list of 20 most related best books with intro, in this format:   {book name}/{intro}.  Final output are in the following format:     - item 1     - item 2     - item 3 ......"
$ws.Range("F7").Value = "任务/请精心设计20个主要任务和主要完成的步骤。回答样式:{任务介绍}-{详细介绍主要完成步骤和方法}  Final output are in the following format:     - 段落 1     - 段落 2     - 段落 3 ......"
$ws.Range("F8").Value = "建议/请精心提供20个最好的建议和具体细节。回答样式:{建议}-{建议的具体组成部分和细节内容}  Final output are in the following format:     - 段落 1     - 段落 2     - 段落 3 ......"

$ws.Range("F9").Select()
